$d = $word.ActiveDocument

function Replace-UniqueText($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
    # Re-find the just-replaced text and toggle a formatting property off/on
    # so the run boundary against its neighbours is preserved (the host merges
    # touching runs that share identical formatting after every text edit).
    $r2 = $d.Content
    $r2.Find.Execute($new, $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
    $r2.Bold = 1
    $r2.Bold = 0
}

function Insert-RunsAfter($anchor, $newTexts) {
    # Locate the single, unambiguous anchor run (a full sentence, never a bare
    # '.') and append each new run's text right after it, in order, each one
    # immediately pinned (Bold off/on) so it doesn't get folded back into its
    # neighbours by the same-formatting run coalescing the host performs after
    # every text edit.
    $r = $d.Content
    $r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null
    $r.Collapse(0)
    $insertAt = $r.Start
    foreach ($t in $newTexts) {
        $ins = $d.Range($insertAt, $insertAt)
        $ins.InsertAfter($t)
        $newEnd = $insertAt + $t.Length
        $pin = $d.Range($insertAt, $newEnd)
        $pin.Bold = 1
        $pin.Bold = 0
        $insertAt = $newEnd
    }
}

# --- Title / author text swaps ---
Replace-UniqueText "The Enchanted Realm of Digital Art" "Democracy: A Collective Choice for a Harmonious Society"
Replace-UniqueText "Isabella Summers" "Alexis Brown"

# --- Author email line: isabella.summers@artech.edu (5 runs) -> at (1 run) ---
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$r3.MoveEnd(1, -1)  # exclude the paragraph mark
$r3.Text = "at"

# --- Body paragraph sentence swaps ---
Replace-UniqueText "In the realm of art, where creativity knows no bounds, a new era has dawned, inviting us into the ethereal realm of digital art" "In the intricate tapestry of human civilization, democracy stands as a beacon of hope, a testament to our innate desire for self-governance and collective progress"
Replace-UniqueText " This remarkable form of artistic expression captivates our minds with its boundless possibilities and challenges our traditional notions of what art can be" " A system born from the collective will of a people, it weaves together the individual threads of aspiration, offering a symphony of voices harmonized in the pursuit of a shared destiny"
Replace-UniqueText " As we venture into the enchanting tapestry of digital art, let us unveil the secrets that lie within its pixels and explore the transformative power it holds in shaping the future of art and technology" " Democracy, with its foundations rooted in the principles of equality, liberty, and fraternity, invites us to participate in the intricate dance of decision-making, challenging us to contemplate the delicate balance between personal autonomy and collective welfare"
Replace-UniqueText "Digital art is an immersive symphony of colors, forms, and textures woven together through the magic of digital tools" "This intricate system, however, is not without its challenges"
Replace-UniqueText " It transcends the limitations of physical mediums, empowering artists to conjure worlds beyond imagination, unfettered by the constraints of materials and space" " The complexities of human nature often weave a web of conflicting interests, leading to debates, disagreements, and even dissent"
Replace-UniqueText " This extraordinary art form has irrevocably altered the landscape of artistic expression, offering a kaleidoscope of possibilities that were once unimaginable" " Yet, within this crucible of diverse perspectives, democracy offers us the profound opportunity to engage in thoughtful discourse, to listen, to understand, and to compromise for the greater good"
Replace-UniqueText "From the ethereal strokes of a digital brush to the mesmerizing fluidity of animated creations, digital art possesses an unrivaled versatility" "Furthermore, democracy empowers us with the solemn responsibility of holding our elected representatives accountable for their actions"
Replace-UniqueText " Artists can now weave intricate patterns and textures, manipulate perspectives, and create an illusion of depth and dimension that is simply not possible with traditional methods" " We, the people, serve as the ultimate custodians of our shared destiny, wielding the power of the ballot box as a potent instrument of accountability"
Replace-UniqueText " This boundless freedom invites experimentation, allowing artists to push the boundaries of their creativity and explore new frontiers of self-expression" " Through this process, we ensure that those who wield authority are ever mindful of the trust we have bestowed upon them and remain steadfast in their commitment to the principles of good governance"
Replace-UniqueText "Digital art has woven its way into the fabric of modern artistic expression, transforming the way we create, appreciate, and interact with art" "In this essay, we have explored the profound significance of democracy, highlighting its role as a collective choice for a harmonious society"
Replace-UniqueText " With its boundless versatility and transformative potential, digital art invites artists and audiences alike to embark on an exhilarating journey through the realm of creativity and imagination" " We have contemplated the intricate balance between individual autonomy and collective welfare, acknowledging the "
Replace-UniqueText " It has opened up new avenues of artistic expression, blurring the boundaries between reality and fantasy, tradition and innovation" " Yet, amidst these challenges, democracy offers us the opportunity to engage in thoughtful discourse, to listen, understand, and compromise for the greater good"
Replace-UniqueText " As technology continues to evolve, the future of digital art is limitless, promising an ever-expanding realm of possibilities where art and technology harmoniously converge" " Moreover, it empowers us with the solemn responsibility of holding our elected representatives accountable, thus ensuring that the decisions made reflect the values and aspirations of all"

# --- New sentences/runs inserted into the body ---
Insert-RunsAfter " Yet, within this crucible of diverse perspectives, democracy offers us the profound opportunity to engage in thoughtful discourse, to listen, to understand, and to compromise for the greater good" @(".", " It is within this arena of civilized dialogue that the true essence of democracy thrives, fostering a spirit of unity amidst diversity, and ensuring that the decisions we make as a collective reflect the values and aspirations of all")
Insert-RunsAfter " We have contemplated the intricate balance between individual autonomy and collective welfare, acknowledging the " @("challenges that arise from the complexities of human nature")

# --- Trailing empty paragraph added at the end of the document ---
$d.Paragraphs.Add() | Out-Null

